$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1117.5769
$ws.Range("I53").Value = 1021.13336
$ws.Range("J53").Value = 1249.091
$ws.Range("K53").Value = 1021.13336
$ws.Range("L53").Value = 1249.091
$ws.Range("M53").Value = -384.13336
$ws.Range("N53").Value = -2523.091

$ws.Range("H138").Value = 2844.5
$ws.Range("I138").Value = 1956.4348
$ws.Range("K138").Value = 5869.3044
$ws.Range("M138").Value = -729.3044

$ws.Range("H140").Value = 214998.86
$ws.Range("J140").Value = 236665.33
$ws.Range("L140").Value = 236665.33
$ws.Range("N140").Value = -247025.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8683039
$ws.Range("I32").Value = 4630768.5
$ws.Range("K32").Value = 4630768.5
$ws.Range("M32").Value = -4630481.5

$ws.Range("H74").Value = 2789.6875
$ws.Range("I74").Value = 2331.1428
$ws.Range("K74").Value = 2331.1428
$ws.Range("M74").Value = -1457.1428

$ws.Range("H77").Value = 2789.6875
$ws.Range("I77").Value = 2331.1428
$ws.Range("K77").Value = 11655.714
$ws.Range("M77").Value = -7287.714

$ws.Range("H88").Value = 2925.6
$ws.Range("J88").Value = 4191.6665
$ws.Range("L88").Value = 4191.6665
$ws.Range("N88").Value = -5003.6665

$ws.Range("H91").Value = 2925.6
$ws.Range("J91").Value = 4191.6665
$ws.Range("L91").Value = 4191.6665
$ws.Range("N91").Value = -6999.6665

$ws.Range("H121").Value = 50240.668
$ws.Range("J121").Value = 50240.668
$ws.Range("L121").Value = 50240.668
$ws.Range("N121").Value = -53734.668

$ws.Range("H128").Value = 50995
$ws.Range("J128").Value = 50995
$ws.Range("L128").Value = 50995
$ws.Range("N128").Value = -60955

$ws.Range("H132").Value = 2430.1072
$ws.Range("I132").Value = 2126.093
$ws.Range("K132").Value = 6378.279
$ws.Range("M132").Value = -3848.279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 23666.166
$ws.Range("I56").Value = 21997
$ws.Range("K56").Value = 21997
$ws.Range("M56").Value = -21258

$ws.Range("H86").Value = 1535.75
$ws.Range("J86").Value = 1303.1538
$ws.Range("L86").Value = 1303.1538
$ws.Range("N86").Value = -3549.1538

$ws.Range("H89").Value = 1535.75
$ws.Range("J89").Value = 1303.1538
$ws.Range("L89").Value = 6515.769
$ws.Range("N89").Value = -17747.769

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H99").Value = 3035.3333
$ws.Range("I99").Value = 2974.8333
$ws.Range("J99").Value = 3277.3333
$ws.Range("K99").Value = 2974.8333
$ws.Range("L99").Value = 3277.3333
$ws.Range("M99").Value = -1476.8333
$ws.Range("N99").Value = -6273.3333

$ws.Range("H107").Value = 1684.3846
$ws.Range("I107").Value = 1533.0834
$ws.Range("K107").Value = 1533.0834
$ws.Range("M107").Value = 386.9166

$ws.Range("H111").Value = 93489.5
$ws.Range("J111").Value = 109979
$ws.Range("L111").Value = 109979
$ws.Range("N111").Value = -118159

$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H134").Value = 11566729
$ws.Range("I134").Value = 2464557
$ws.Range("K134").Value = 7393671
$ws.Range("M134").Value = -7391136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 535.13336
$ws.Range("I7").Value = 501.92856
$ws.Range("K7").Value = 501.92856
$ws.Range("M7").Value = -388.92856

$ws.Range("H16").Value = 1284.6316
$ws.Range("I16").Value = 1236.1666
$ws.Range("J16").Value = 1367.7142
$ws.Range("K16").Value = 1236.1666
$ws.Range("L16").Value = 1367.7142
$ws.Range("M16").Value = -949.1666
$ws.Range("N16").Value = -1941.7142

$ws.Range("H18").Value = 40573
$ws.Range("J18").Value = 40573
$ws.Range("L18").Value = 40573
$ws.Range("N18").Value = -41033

$ws.Range("H31").Value = 5207.5
$ws.Range("I31").Value = 2973.6
$ws.Range("J31").Value = 6066.6924
$ws.Range("K31").Value = 2973.6
$ws.Range("L31").Value = 6066.6924
$ws.Range("M31").Value = -2678.6
$ws.Range("N31").Value = -6656.6924

$ws.Range("H34").Value = 5207.5
$ws.Range("I34").Value = 2973.6
$ws.Range("J34").Value = 6066.6924
$ws.Range("K34").Value = 2973.6
$ws.Range("L34").Value = 6066.6924
$ws.Range("M34").Value = -2771.6
$ws.Range("N34").Value = -6470.6924

$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 10000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -11186

$ws.Range("H50").Value = 36166.668
$ws.Range("J50").Value = 48500
$ws.Range("L50").Value = 48500
$ws.Range("N50").Value = -49750

$ws.Range("H58").Value = 2813.0356
$ws.Range("I58").Value = 2252.6316
$ws.Range("K58").Value = 2252.6316
$ws.Range("M58").Value = -2049.6316

$ws.Range("H75").Value = 104640
$ws.Range("J75").Value = 104640
$ws.Range("L75").Value = 104640
$ws.Range("N75").Value = -106636

$ws.Range("H78").Value = 104640
$ws.Range("J78").Value = 104640
$ws.Range("L78").Value = 313920
$ws.Range("N78").Value = -323904

$ws.Range("H99").Value = 1495.75
$ws.Range("I99").Value = 1495.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1495.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 2.25
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 108995
$ws.Range("J100").Value = 108995
$ws.Range("L100").Value = 108995
$ws.Range("N100").Value = -111159

$ws.Range("H113").Value = 1284.6316
$ws.Range("I113").Value = 1236.1666
$ws.Range("J113").Value = 1367.7142
$ws.Range("K113").Value = 1236.1666
$ws.Range("L113").Value = 1367.7142
$ws.Range("M113").Value = 933.8334
$ws.Range("N113").Value = -5707.7142

$ws.Range("H126").Value = 1495.75
$ws.Range("I126").Value = 1495.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4487.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2017.25
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 2813.0356
$ws.Range("I136").Value = 2252.6316
$ws.Range("K136").Value = 6757.8948
$ws.Range("M136").Value = -4207.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 817.44446
$ws.Range("J122").Value = 720.6667
$ws.Range("L122").Value = 6486.0003
$ws.Range("N122").Value = -11386.0003

$ws.Range("H128").Value = 78265.336
$ws.Range("I128").Value = 78265.336
$ws.Range("K128").Value = 234796.008
$ws.Range("M128").Value = -229816.008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 90978
$ws.Range("J51").Value = 90978
$ws.Range("L51").Value = 90978
$ws.Range("N51").Value = -91996

$ws.Range("H97").Value = 837.5925999999999
$ws.Range("I97").Value = 477.04544
$ws.Range("K97").Value = 477.04544
$ws.Range("M97").Value = 18.95456000000001

$ws.Range("H128").Value = 134719
$ws.Range("J128").Value = 134719
$ws.Range("L128").Value = 134719
$ws.Range("N128").Value = -144679

$ws.Range("H132").Value = 1829.45
$ws.Range("I132").Value = 1327.9412
$ws.Range("K132").Value = 3983.8236
$ws.Range("M132").Value = -1453.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17988.572
$ws.Range("I7").Value = 17922.562
$ws.Range("K7").Value = 17922.562
$ws.Range("M7").Value = -17810.562

$ws.Range("H46").Value = 2836.5
$ws.Range("I46").Value = 1319.6
$ws.Range("K46").Value = 1319.6
$ws.Range("M46").Value = -1131.6

$ws.Range("H82").Value = 2421.348
$ws.Range("I82").Value = 1873.3334
$ws.Range("J82").Value = 3019.182
$ws.Range("K82").Value = 1873.3334
$ws.Range("L82").Value = 3019.182
$ws.Range("M82").Value = -1512.3334
$ws.Range("N82").Value = -3741.182

$ws.Range("H85").Value = 2421.348
$ws.Range("I85").Value = 1873.3334
$ws.Range("J85").Value = 3019.182
$ws.Range("K85").Value = 1873.3334
$ws.Range("L85").Value = 3019.182
$ws.Range("M85").Value = -625.3334
$ws.Range("N85").Value = -5515.182

$ws.Range("H100").Value = 5401
$ws.Range("I100").Value = 4601.5
$ws.Range("K100").Value = 4601.5
$ws.Range("M100").Value = -4060.5

$ws.Range("H122").Value = 3242.2632
$ws.Range("I122").Value = 3242.2632
$ws.Range("K122").Value = 9726.7896
$ws.Range("M122").Value = -7276.7896

$ws.Range("H126").Value = 17988.572
$ws.Range("I126").Value = 17922.562
$ws.Range("K126").Value = 53767.686
$ws.Range("M126").Value = -51297.686

$ws.Range("H141").Value = 999999
$ws.Range("J141").Value = 999999
$ws.Range("L141").Value = 999999
$ws.Range("N141").Value = -1010359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 46434.5
$ws.Range("J16").Value = 46434.5
$ws.Range("L16").Value = 46434.5
$ws.Range("N16").Value = -47018.5

$ws.Range("H136").Value = 24209.744
$ws.Range("I136").Value = 1790.5385
$ws.Range("K136").Value = 5371.6155
$ws.Range("M136").Value = -2821.6155
